$p = $ppt.ActivePresentation

# --- Change 1: Slide 2 ("Introduction") ---
# "Some FL Techniques: ... Mutation-Based" -> "... Mutation-Based FL"
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(6)
$tr2 = $shape2.TextFrame.TextRange
$run2 = $tr2.Paragraphs(3).Runs(1)
$run2.Text = "Some FL Techniques: Spectrum Based FL, Slice-Based FL, Mutation-Based FL"

# --- Change 2: Slide 5 ("Key Idea") ---
# Merge the two runs of the "Analyze which categories..." bullet into a single run.
$slide5 = $p.Slides.Item(5)
$shape5 = $slide5.Shapes.Item(6)
$tr5 = $shape5.TextFrame.TextRange
$para5 = $tr5.Paragraphs(5)
$run5a = $para5.Runs(1)
$run5b = $para5.Runs(2)
$run5a.Text = "Analyze which categories of faults (based on taxonomy) are most effectively localized, focusing on the traits of faults that successfully identified vs which are not which are not identified. "
$run5b.Text = ""
